# Applies the scheduled market-data refresh to the Leve profit sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns), mirroring
# the upstream API pull that is periodically re-run by the commit's
# scheduled runner. Pure value updates -- no formulas, no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 38296.668
$ws.Range("J44").Value = 38296.668
$ws.Range("L44").Value = 38296.668
$ws.Range("N44").Value = -39220.668
$ws.Range("H62").Value = 9067.799999999999
$ws.Range("I62").Value = 7666.3335
$ws.Range("J62").Value = 11170
$ws.Range("K62").Value = 7666.3335
$ws.Range("L62").Value = 11170
$ws.Range("M62").Value = -7042.3335
$ws.Range("N62").Value = -12418
$ws.Range("H65").Value = 9067.799999999999
$ws.Range("I65").Value = 7666.3335
$ws.Range("J65").Value = 11170
$ws.Range("K65").Value = 38331.6675
$ws.Range("L65").Value = 55850
$ws.Range("M65").Value = -35211.6675
$ws.Range("N65").Value = -62090
$ws.Range("H69").Value = 10015
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 10015
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H86").Value = 83380010
$ws.Range("I86").Value = 66679348
$ws.Range("J86").Value = 104255830
$ws.Range("K86").Value = 66679348
$ws.Range("L86").Value = 104255830
$ws.Range("M86").Value = -66678225
$ws.Range("N86").Value = -104258076
$ws.Range("H89").Value = 83380010
$ws.Range("I89").Value = 66679348
$ws.Range("J89").Value = 104255830
$ws.Range("K89").Value = 333396740
$ws.Range("L89").Value = 521279150
$ws.Range("M89").Value = -333391124
$ws.Range("N89").Value = -521290382
$ws.Range("H92").Value = 16129834
$ws.Range("I92").Value = 20000844
$ws.Range("J92").Value = 624.5
$ws.Range("K92").Value = 20000844
$ws.Range("L92").Value = 624.5
$ws.Range("M92").Value = -19999596
$ws.Range("N92").Value = -3120.5
$ws.Range("H98").Value = 1047.9286
$ws.Range("I98").Value = 472.66666
$ws.Range("J98").Value = 4499.5
$ws.Range("K98").Value = 472.66666
$ws.Range("L98").Value = 4499.5
$ws.Range("M98").Value = 1025.33334
$ws.Range("N98").Value = -7495.5
$ws.Range("H100").Value = 2239.0557
$ws.Range("I100").Value = 1572.4783
$ws.Range("J100").Value = 3418.3845
$ws.Range("K100").Value = 1572.4783
$ws.Range("L100").Value = 3418.3845
$ws.Range("M100").Value = -1031.4783
$ws.Range("N100").Value = -4500.3845
$ws.Range("H113").Value = 10472.417
$ws.Range("I113").Value = 14413
$ws.Range("J113").Value = 2591.25
$ws.Range("K113").Value = 14413
$ws.Range("L113").Value = 2591.25
$ws.Range("M113").Value = -11159
$ws.Range("N113").Value = -9099.25
$ws.Range("H122").Value = 1047.9286
$ws.Range("I122").Value = 472.66666
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 1417.99998
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = 1032.00002
$ws.Range("N122").Value = -18398.5
$ws.Range("H137").Value = 49327.39
$ws.Range("I137").Value = 77918.27
$ws.Range("K137").Value = 233754.81
$ws.Range("M137").Value = -231204.81
$ws.Range("H138").Value = 2003.68
$ws.Range("I138").Value = 1081.9429
$ws.Range("K138").Value = 3245.8287
$ws.Range("M138").Value = 1894.1713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23317736
$ws.Range("I32").Value = 32014024
$ws.Range("K32").Value = 32014024
$ws.Range("M32").Value = -32013737
$ws.Range("H61").Value = 4270.478
$ws.Range("I61").Value = 4541.7334
$ws.Range("K61").Value = 4541.7334
$ws.Range("M61").Value = -4329.7334
$ws.Range("H74").Value = 2531.4666
$ws.Range("I74").Value = 2369.8206
$ws.Range("K74").Value = 2369.8206
$ws.Range("M74").Value = -1495.8206
$ws.Range("H77").Value = 2531.4666
$ws.Range("I77").Value = 2369.8206
$ws.Range("K77").Value = 11849.103
$ws.Range("M77").Value = -7481.102999999999
$ws.Range("H92").Value = 67929.664
$ws.Range("J92").Value = 67929.664
$ws.Range("L92").Value = 67929.664
$ws.Range("N92").Value = -72921.664
$ws.Range("H136").Value = 4270.478
$ws.Range("I136").Value = 4541.7334
$ws.Range("K136").Value = 13625.2002
$ws.Range("M136").Value = -11075.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2991.5454
$ws.Range("I99").Value = 2251.3333
$ws.Range("K99").Value = 2251.3333
$ws.Range("M99").Value = -753.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 85.25
$ws.Range("I7").Value = 33.666668
$ws.Range("K7").Value = 33.666668
$ws.Range("M7").Value = 79.333332
$ws.Range("H31").Value = 5809.4136
$ws.Range("I31").Value = 2558.5
$ws.Range("K31").Value = 2558.5
$ws.Range("M31").Value = -2263.5
$ws.Range("H34").Value = 5809.4136
$ws.Range("I34").Value = 2558.5
$ws.Range("K34").Value = 2558.5
$ws.Range("M34").Value = -2356.5
$ws.Range("H58").Value = 2694.325
$ws.Range("I58").Value = 2635.205
$ws.Range("K58").Value = 2635.205
$ws.Range("M58").Value = -2432.205
$ws.Range("H86").Value = 41334.39
$ws.Range("I86").Value = 40332.332
$ws.Range("K86").Value = 40332.332
$ws.Range("M86").Value = -39209.332
$ws.Range("H89").Value = 41334.39
$ws.Range("I89").Value = 40332.332
$ws.Range("K89").Value = 201661.66
$ws.Range("M89").Value = -196045.66
$ws.Range("H132").Value = 5005.5
$ws.Range("I132").Value = 4996.294
$ws.Range("J132").Value = 5027.857
$ws.Range("K132").Value = 14988.882
$ws.Range("L132").Value = 15083.571
$ws.Range("M132").Value = -12458.882
$ws.Range("N132").Value = -20143.571
$ws.Range("H134").Value = 2491.6155
$ws.Range("J134").Value = 3333
$ws.Range("L134").Value = 9999
$ws.Range("N134").Value = -15069
$ws.Range("H136").Value = 2694.325
$ws.Range("I136").Value = 2635.205
$ws.Range("K136").Value = 7905.615
$ws.Range("M136").Value = -5355.615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2.6
$ws.Range("J12").Value = 2.6
$ws.Range("L12").Value = 7.800000000000001
$ws.Range("N12").Value = -353.8
$ws.Range("H129").Value = 2516.2222
$ws.Range("I129").Value = 795
$ws.Range("J129").Value = 3008
$ws.Range("K129").Value = 2385
$ws.Range("L129").Value = 9024
$ws.Range("M129").Value = 2615
$ws.Range("N129").Value = -19024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6500284
$ws.Range("J3").Value = 10000000
$ws.Range("L3").Value = 10000000
$ws.Range("N3").Value = -10000232
$ws.Range("H49").Value = 20000000
$ws.Range("I49").Value = 20000000
$ws.Range("K49").Value = 20000000
$ws.Range("M49").Value = -19999816
$ws.Range("H97").Value = 1103.25
$ws.Range("I97").Value = 709.7143
$ws.Range("K97").Value = 709.7143
$ws.Range("M97").Value = -213.7143
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H120").Value = 88278.336
$ws.Range("J120").Value = 88278.336
$ws.Range("L120").Value = 88278.336
$ws.Range("N120").Value = -97954.336
$ws.Range("H132").Value = 4110
$ws.Range("I132").Value = 4118.5415
$ws.Range("J132").Value = 4091.3635
$ws.Range("K132").Value = 12355.6245
$ws.Range("L132").Value = 12274.0905
$ws.Range("M132").Value = -9825.624500000002
$ws.Range("N132").Value = -17334.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 30025
$ws.Range("I42").Value = 30025
$ws.Range("K42").Value = 30025
$ws.Range("M42").Value = -29462
$ws.Range("H49").Value = 30025
$ws.Range("I49").Value = 30025
$ws.Range("K49").Value = 30025
$ws.Range("M49").Value = -29878
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H68").Value = 7532.467
$ws.Range("J68").Value = 7714.7144
$ws.Range("L68").Value = 7714.7144
$ws.Range("N68").Value = -9212.714400000001
$ws.Range("H71").Value = 7532.467
$ws.Range("J71").Value = 7714.7144
$ws.Range("L71").Value = 38573.572
$ws.Range("N71").Value = -46061.572
$ws.Range("H122").Value = 11495.6
$ws.Range("I122").Value = 12013.182
$ws.Range("J122").Value = 10072.25
$ws.Range("K122").Value = 36039.546
$ws.Range("L122").Value = 30216.75
$ws.Range("M122").Value = -33589.546
$ws.Range("N122").Value = -35116.75
$ws.Range("H132").Value = 554932.4399999999
$ws.Range("I132").Value = 921249
$ws.Range("J132").Value = 5457.625
$ws.Range("K132").Value = 2763747
$ws.Range("L132").Value = 16372.875
$ws.Range("M132").Value = -2761217
$ws.Range("N132").Value = -21432.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 62291.234
$ws.Range("I132").Value = 69896.39999999999
$ws.Range("K132").Value = 209689.2
$ws.Range("M132").Value = -207159.2
$ws.Range("H136").Value = 40741.27
$ws.Range("I136").Value = 1456.5555
$ws.Range("K136").Value = 4369.666499999999
$ws.Range("M136").Value = -1819.666499999999
